# Auto-generated script applying scheduled-runner market data refresh
# to the Sagittarius_Profits workbook's per-job sheets (ALC/ARM/BSM/CRP/CUL/GSM/LTW/WVR).
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H2").Value = 4030.1667
$ws.Range("I2").Value = 400
$ws.Range("J2").Value = 7660.3335
$ws.Range("K2").Value = 400
$ws.Range("L2").Value = 7660.3335
$ws.Range("M2").Value = -287
$ws.Range("N2").Value = -7886.3335
$ws.Range("H28").Value = 1942.75
$ws.Range("I28").Value = 2323.8333
$ws.Range("J28").Value = 799.5
$ws.Range("K28").Value = 2323.8333
$ws.Range("L28").Value = 799.5
$ws.Range("M28").Value = -1838.8333
$ws.Range("N28").ClearContents()
$ws.Range("H64").Value = 4497.25
$ws.Range("I64").Value = 3489
$ws.Range("J64").Value = 4833.3335
$ws.Range("K64").Value = 3489
$ws.Range("L64").Value = 4833.3335
$ws.Range("M64").Value = -3241
$ws.Range("N64").Value = -5329.3335
$ws.Range("H67").Value = 4497.25
$ws.Range("I67").Value = 3489
$ws.Range("J67").Value = 4833.3335
$ws.Range("K67").Value = 3489
$ws.Range("L67").Value = 4833.3335
$ws.Range("M67").Value = -2631
$ws.Range("N67").Value = -6549.3335
$ws.Range("H92").Value = 853.4286
$ws.Range("I92").Value = 631.2727
$ws.Range("J92").Value = 1668
$ws.Range("K92").Value = 631.2727
$ws.Range("L92").Value = 1668
$ws.Range("M92").Value = 616.7273
$ws.Range("N92").Value = -4164
$ws.Range("H98").Value = 1991.75
$ws.Range("I98").Value = 1991.75
$ws.Range("J98").Value = 0
$ws.Range("K98").Value = 1991.75
$ws.Range("L98").Value = 0
$ws.Range("M98").Value = -493.75
$ws.Range("H111").Value = 8252
$ws.Range("I111").Value = 8483.727999999999
$ws.Range("J111").Value = 7402.3335
$ws.Range("K111").Value = 25451.184
$ws.Range("L111").Value = 22207.0005
$ws.Range("M111").Value = -22384.184
$ws.Range("N111").Value = -28341.0005
$ws.Range("H122").Value = 1991.75
$ws.Range("I122").Value = 1991.75
$ws.Range("J122").Value = 0
$ws.Range("K122").Value = 5975.25
$ws.Range("L122").Value = 0
$ws.Range("M122").Value = -3525.25
$ws.Range("H132").Value = 2470.5
$ws.Range("I132").Value = 2632.0833
$ws.Range("J132").Value = 1985.75
$ws.Range("K132").Value = 7896.249899999999
$ws.Range("L132").Value = 5957.25
$ws.Range("M132").Value = -5366.249899999999
$ws.Range("N132").Value = -11017.25
$ws.Range("H135").Value = 2331.6875
$ws.Range("I135").Value = 2287.1333
$ws.Range("J135").Value = 3000
$ws.Range("K135").Value = 20584.1997
$ws.Range("L135").Value = 27000
$ws.Range("M135").Value = -18049.1997
$ws.Range("N135").ClearContents()
$ws.Range("H137").Value = 2225
$ws.Range("I137").Value = 1950
$ws.Range("J137").Value = 2500
$ws.Range("K137").Value = 5850
$ws.Range("L137").Value = 7500
$ws.Range("M137").Value = -3300
$ws.Range("N137").ClearContents()
$ws.Range("H138").Value = 5782.722
$ws.Range("I138").Value = 5999
$ws.Range("J138").Value = 5770
$ws.Range("K138").Value = 17997
$ws.Range("L138").Value = 17310
$ws.Range("M138").Value = -12857
$ws.Range("N138").Value = -27590

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 5137815.5
$ws.Range("I32").Value = 5004802
$ws.Range("J32").Value = 7000000
$ws.Range("K32").Value = 5004802
$ws.Range("L32").Value = 7000000
$ws.Range("M32").Value = -5004515
$ws.Range("N32").ClearContents()
$ws.Range("H61").Value = 1666.3334
$ws.Range("I61").Value = 1666.3334
$ws.Range("J61").Value = 0
$ws.Range("K61").Value = 1666.3334
$ws.Range("L61").Value = 0
$ws.Range("M61").Value = -1454.3334
$ws.Range("H68").Value = 0
$ws.Range("I68").Value = 0
$ws.Range("J68").Value = 0
$ws.Range("K68").Value = 0
$ws.Range("L68").Value = 0
$ws.Range("N68").ClearContents()
$ws.Range("H71").Value = 0
$ws.Range("I71").Value = 0
$ws.Range("J71").Value = 0
$ws.Range("K71").Value = 0
$ws.Range("L71").Value = 0
$ws.Range("N71").ClearContents()
$ws.Range("H92").Value = 56745
$ws.Range("I92").Value = 0
$ws.Range("J92").Value = 56745
$ws.Range("K92").Value = 0
$ws.Range("L92").Value = 56745
$ws.Range("N92").Value = -61737
$ws.Range("H98").Value = 0
$ws.Range("I98").Value = 0
$ws.Range("J98").Value = 0
$ws.Range("K98").Value = 0
$ws.Range("L98").Value = 0
$ws.Range("N98").ClearContents()
$ws.Range("H136").Value = 1666.3334
$ws.Range("I136").Value = 1666.3334
$ws.Range("J136").Value = 0
$ws.Range("K136").Value = 4999.0002
$ws.Range("L136").Value = 0
$ws.Range("M136").Value = -2449.0002

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H94").Value = 2076.6667
$ws.Range("I94").Value = 1909.5
$ws.Range("J94").Value = 2411
$ws.Range("K94").Value = 1909.5
$ws.Range("L94").Value = 2411
$ws.Range("M94").Value = -1458.5
$ws.Range("N94").Value = -3313
$ws.Range("H107").Value = 2056
$ws.Range("I107").Value = 574.8333
$ws.Range("J107").Value = 6499.5
$ws.Range("K107").Value = 574.8333
$ws.Range("L107").Value = 6499.5
$ws.Range("M107").Value = 1345.1667
$ws.Range("N107").ClearContents()

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H7").Value = 336.55554
$ws.Range("I7").Value = 140.875
$ws.Range("J7").Value = 1902
$ws.Range("K7").Value = 140.875
$ws.Range("L7").Value = 1902
$ws.Range("M7").Value = -27.875
$ws.Range("N7").ClearContents()
$ws.Range("H22").Value = 2333
$ws.Range("I22").Value = 2333
$ws.Range("J22").Value = 0
$ws.Range("K22").Value = 2333
$ws.Range("L22").Value = 0
$ws.Range("M22").Value = -1983
$ws.Range("N22").ClearContents()
$ws.Range("H58").Value = 3996.3333
$ws.Range("I58").Value = 2106
$ws.Range("J58").Value = 7777
$ws.Range("K58").Value = 2106
$ws.Range("L58").Value = 7777
$ws.Range("M58").Value = -1903
$ws.Range("N58").Value = -8183
$ws.Range("H99").Value = 864.6842
$ws.Range("I99").Value = 834.53845
$ws.Range("J99").Value = 930
$ws.Range("K99").Value = 834.53845
$ws.Range("L99").Value = 930
$ws.Range("M99").Value = 663.46155
$ws.Range("N99").Value = -3926
$ws.Range("H126").Value = 864.6842
$ws.Range("I126").Value = 834.53845
$ws.Range("J126").Value = 930
$ws.Range("K126").Value = 2503.61535
$ws.Range("L126").Value = 2790
$ws.Range("M126").Value = -33.61535000000003
$ws.Range("N126").Value = -7730
$ws.Range("H134").Value = 2730.7778
$ws.Range("I134").Value = 2515.4
$ws.Range("J134").Value = 3000
$ws.Range("K134").Value = 7546.200000000001
$ws.Range("L134").Value = 9000
$ws.Range("M134").Value = -5011.200000000001
$ws.Range("N134").ClearContents()
$ws.Range("H136").Value = 3996.3333
$ws.Range("I136").Value = 2106
$ws.Range("J136").Value = 7777
$ws.Range("K136").Value = 6318
$ws.Range("L136").Value = 23331
$ws.Range("M136").Value = -3768
$ws.Range("N136").Value = -28431

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H68").Value = 4682.846
$ws.Range("I68").Value = 3500
$ws.Range("J68").Value = 4781.4165
$ws.Range("K68").Value = 10500
$ws.Range("L68").Value = 14344.2495
$ws.Range("M68").Value = -9689
$ws.Range("N68").Value = -15966.2495
$ws.Range("H71").Value = 4682.846
$ws.Range("I71").Value = 3500
$ws.Range("J71").Value = 4781.4165
$ws.Range("K71").Value = 31500
$ws.Range("L71").Value = 43032.7485
$ws.Range("M71").Value = -27444
$ws.Range("N71").Value = -51144.7485
$ws.Range("H132").Value = 1097.4
$ws.Range("I132").Value = 1497.3334
$ws.Range("J132").Value = 497.5
$ws.Range("K132").Value = 13476.0006
$ws.Range("L132").Value = 4477.5
$ws.Range("M132").Value = -10946.0006
$ws.Range("N132").Value = -9537.5
$ws.Range("H137").Value = 2810.2856
$ws.Range("I137").Value = 2392.5
$ws.Range("J137").Value = 2977.4
$ws.Range("K137").Value = 7177.5
$ws.Range("L137").Value = 8932.200000000001
$ws.Range("M137").Value = -2077.5
$ws.Range("N137").Value = -19132.2
$ws.Range("H138").Value = 1117
$ws.Range("I138").Value = 1117
$ws.Range("J138").Value = 0
$ws.Range("K138").Value = 3351
$ws.Range("L138").Value = 0
$ws.Range("M138").Value = 1789

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H2").Value = 183.66667
$ws.Range("I2").Value = 100.5
$ws.Range("J2").Value = 350
$ws.Range("K2").Value = 100.5
$ws.Range("L2").Value = 350
$ws.Range("M2").Value = 12.5
$ws.Range("N2").Value = -576
$ws.Range("H43").Value = 30000
$ws.Range("I43").Value = 0
$ws.Range("J43").Value = 30000
$ws.Range("K43").Value = 0
$ws.Range("L43").Value = 30000
$ws.Range("M43").ClearContents()
$ws.Range("N43").Value = -30302
$ws.Range("H80").Value = 27124.25
$ws.Range("I80").Value = 3250
$ws.Range("J80").Value = 50998.5
$ws.Range("K80").Value = 3250
$ws.Range("L80").Value = 50998.5
$ws.Range("M80").Value = -2252
$ws.Range("N80").Value = -52994.5
$ws.Range("H83").Value = 27124.25
$ws.Range("I83").Value = 3250
$ws.Range("J83").Value = 50998.5
$ws.Range("K83").Value = 16250
$ws.Range("L83").Value = 254992.5
$ws.Range("M83").Value = -11258
$ws.Range("N83").Value = -264976.5
$ws.Range("H107").Value = 1749.6666
$ws.Range("I107").Value = 89
$ws.Range("J107").Value = 3825.5
$ws.Range("K107").Value = 89
$ws.Range("L107").Value = 3825.5
$ws.Range("M107").Value = 1831
$ws.Range("N107").ClearContents()
$ws.Range("H122").Value = 3467
$ws.Range("I122").Value = 3467
$ws.Range("J122").Value = 0
$ws.Range("K122").Value = 10401
$ws.Range("L122").Value = 0
$ws.Range("M122").Value = -7951
$ws.Range("N122").ClearContents()

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 7243.375
$ws.Range("I7").Value = 5450
$ws.Range("J7").Value = 7499.5713
$ws.Range("K7").Value = 5450
$ws.Range("L7").Value = 7499.5713
$ws.Range("M7").Value = -5338
$ws.Range("N7").ClearContents()
$ws.Range("H22").Value = 2791.6667
$ws.Range("I22").Value = 2443.5
$ws.Range("J22").Value = 3488
$ws.Range("K22").Value = 2443.5
$ws.Range("L22").Value = 3488
$ws.Range("M22").Value = -2148.5
$ws.Range("N22").Value = -4078
$ws.Range("H27").Value = 2791.6667
$ws.Range("I27").Value = 2443.5
$ws.Range("J27").Value = 3488
$ws.Range("K27").Value = 2443.5
$ws.Range("L27").Value = 3488
$ws.Range("M27").Value = -2336.5
$ws.Range("N27").Value = -3702
$ws.Range("H98").Value = 32495
$ws.Range("I98").Value = 0
$ws.Range("J98").Value = 32495
$ws.Range("K98").Value = 0
$ws.Range("L98").Value = 32495
$ws.Range("N98").Value = -38485
$ws.Range("H126").Value = 7243.375
$ws.Range("I126").Value = 5450
$ws.Range("J126").Value = 7499.5713
$ws.Range("K126").Value = 16350
$ws.Range("L126").Value = 22498.7139
$ws.Range("M126").Value = -13880
$ws.Range("N126").ClearContents()
$ws.Range("H132").Value = 3027.4546
$ws.Range("I132").Value = 2861.5557
$ws.Range("J132").Value = 3774
$ws.Range("K132").Value = 8584.667099999999
$ws.Range("L132").Value = 11322
$ws.Range("M132").Value = -6054.667099999999
$ws.Range("N132").ClearContents()
$ws.Range("H136").Value = 1400
$ws.Range("I136").Value = 1400
$ws.Range("J136").Value = 0
$ws.Range("K136").Value = 4200
$ws.Range("L136").Value = 0
$ws.Range("M136").Value = -1650
$ws.Range("N136").ClearContents()

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H107").Value = 1737.7646
$ws.Range("I107").Value = 1650.2727
$ws.Range("J107").Value = 1898.1666
$ws.Range("K107").Value = 4950.8181
$ws.Range("L107").Value = 5694.4998
$ws.Range("M107").Value = -3030.8181
$ws.Range("N107").Value = -9534.4998
$ws.Range("H136").Value = 2875.6191
$ws.Range("I136").Value = 2827.5
$ws.Range("J136").Value = 3164.3333
$ws.Range("K136").Value = 8482.5
$ws.Range("L136").Value = 9492.999899999999
$ws.Range("M136").Value = -5932.5
$ws.Range("N136").ClearContents()
